$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3564250
$ws.Range("J17").Value = 3564250
$ws.Range("L17").Value = 10692750
$ws.Range("N17").Value = -10693086
$ws.Range("H33").Value = 12357489
$ws.Range("I33").Value = 5117658
$ws.Range("K33").Value = 5117658
$ws.Range("M33").Value = -5117429
$ws.Range("H45").Value = 24471.666
$ws.Range("J45").Value = 24471.666
$ws.Range("L45").Value = 73414.998
$ws.Range("N45").Value = -73798.998
$ws.Range("H51").Value = 69445740
$ws.Range("J51").Value = 125001270
$ws.Range("L51").Value = 125001270
$ws.Range("N51").Value = -125002238
$ws.Range("H76").Value = 3121.6667
$ws.Range("I76").Value = 3087.6
$ws.Range("K76").Value = 3087.6
$ws.Range("M76").Value = -2772.6
$ws.Range("H79").Value = 3121.6667
$ws.Range("I79").Value = 3087.6
$ws.Range("K79").Value = 3087.6
$ws.Range("M79").Value = -1995.6
$ws.Range("H87").Value = 77547.57
$ws.Range("J87").Value = 80710.52
$ws.Range("L87").Value = 80710.52
$ws.Range("N87").Value = -83206.52
$ws.Range("H90").Value = 77547.57
$ws.Range("J90").Value = 80710.52
$ws.Range("L90").Value = 242131.56
$ws.Range("N90").Value = -254611.56
$ws.Range("H134").Value = 118425
$ws.Range("J134").Value = 131495.83
$ws.Range("L134").Value = 131495.83
$ws.Range("N134").Value = -141635.83
$ws.Range("H138").Value = 4594.0137
$ws.Range("I138").Value = 7921.227
$ws.Range("J138").Value = 3158.745
$ws.Range("K138").Value = 23763.681
$ws.Range("L138").Value = 9476.235
$ws.Range("M138").Value = -18623.681
$ws.Range("N138").Value = -19756.235

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2000
$ws.Range("J26").Value = 2000
$ws.Range("L26").Value = 2000
$ws.Range("N26").Value = -2660
$ws.Range("H32").Value = 223593.94
$ws.Range("I32").Value = 265823.9
$ws.Range("K32").Value = 265823.9
$ws.Range("M32").Value = -265536.9
$ws.Range("H74").Value = 1115927.1
$ws.Range("I74").Value = 1974.875
$ws.Range("J74").Value = 1858562
$ws.Range("K74").Value = 1974.875
$ws.Range("L74").Value = 1858562
$ws.Range("M74").Value = -1100.875
$ws.Range("N74").Value = -1860310
$ws.Range("H77").Value = 1115927.1
$ws.Range("I77").Value = 1974.875
$ws.Range("J77").Value = 1858562
$ws.Range("K77").Value = 9874.375
$ws.Range("L77").Value = 9292810
$ws.Range("M77").Value = -5506.375
$ws.Range("N77").Value = -9301546
$ws.Range("H110").Value = 751.8889
$ws.Range("I110").Value = 751.8889
$ws.Range("K110").Value = 751.8889
$ws.Range("M110").Value = 1293.1111
$ws.Range("H132").Value = 2379.3
$ws.Range("J132").Value = 4275
$ws.Range("L132").Value = 12825
$ws.Range("N132").Value = -17885

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2701.7693
$ws.Range("I94").Value = 2283.9092
$ws.Range("K94").Value = 2283.9092
$ws.Range("M94").Value = -1832.9092
$ws.Range("H107").Value = 16784.953
$ws.Range("I107").Value = 18852.588
$ws.Range("K107").Value = 18852.588
$ws.Range("M107").Value = -16932.588

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 461.025
$ws.Range("I22").Value = 488.3143
$ws.Range("J22").Value = 270
$ws.Range("K22").Value = 488.3143
$ws.Range("L22").Value = 270
$ws.Range("M22").Value = -138.3143
$ws.Range("N22").Value = -970
$ws.Range("H31").Value = 2907.879
$ws.Range("I31").Value = 1689.1538
$ws.Range("J31").Value = 3206.8113
$ws.Range("K31").Value = 1689.1538
$ws.Range("L31").Value = 3206.8113
$ws.Range("M31").Value = -1394.1538
$ws.Range("N31").Value = -3796.8113
$ws.Range("H34").Value = 2907.879
$ws.Range("I34").Value = 1689.1538
$ws.Range("J34").Value = 3206.8113
$ws.Range("K34").Value = 1689.1538
$ws.Range("L34").Value = 3206.8113
$ws.Range("M34").Value = -1487.1538
$ws.Range("N34").Value = -3610.8113
$ws.Range("H93").Value = 10333
$ws.Range("I93").Value = 10333
$ws.Range("K93").Value = 10333
$ws.Range("M93").Value = -8461
$ws.Range("H132").Value = 15154629
$ws.Range("I132").Value = 4660
$ws.Range("K132").Value = 13980
$ws.Range("M132").Value = -11450

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1167.7241
$ws.Range("J34").Value = 1245.8518
$ws.Range("L34").Value = 3737.5554
$ws.Range("N34").Value = -3905.5554
$ws.Range("H69").Value = 10387.556
$ws.Range("I69").Value = 18466.666
$ws.Range("J69").Value = 6348
$ws.Range("K69").Value = 55399.99800000001
$ws.Range("L69").Value = 19044
$ws.Range("M69").Value = -54588.99800000001
$ws.Range("N69").Value = -20666
$ws.Range("H72").Value = 10387.556
$ws.Range("I72").Value = 18466.666
$ws.Range("J72").Value = 6348
$ws.Range("K72").Value = 166199.994
$ws.Range("L72").Value = 57132
$ws.Range("M72").Value = -162143.994
$ws.Range("N72").Value = -65244
$ws.Range("H86").Value = 78490.62
$ws.Range("I86").Value = 1273.5
$ws.Range("K86").Value = 3820.5
$ws.Range("M86").Value = -2634.5
$ws.Range("H89").Value = 78490.62
$ws.Range("I89").Value = 1273.5
$ws.Range("K89").Value = 11461.5
$ws.Range("M89").Value = -5533.5
$ws.Range("H99").Value = 100008984
$ws.Range("I99").Value = 333336960
$ws.Range("J99").Value = 11285
$ws.Range("K99").Value = 1000010880
$ws.Range("L99").Value = 33855
$ws.Range("M99").Value = -1000008634
$ws.Range("N99").Value = -38347
$ws.Range("H113").Value = 1071.8889
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1080.875
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 3242.625
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -7582.625
$ws.Range("H131").Value = 9231592
$ws.Range("J131").Value = 7010490
$ws.Range("L131").Value = 21031470
$ws.Range("N131").Value = -21041550

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("H43").Value = 11416.5
$ws.Range("I43").Value = 11700
$ws.Range("J43").Value = 9999
$ws.Range("K43").Value = 11700
$ws.Range("L43").Value = 9999
$ws.Range("M43").Value = -11549
$ws.Range("N43").Value = -10301
$ws.Range("N33").ClearContents()
$ws.Range("N36").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 461796.34
$ws.Range("I23").Value = 22309.572
$ws.Range("J23").Value = 2000000
$ws.Range("K23").Value = 22309.572
$ws.Range("L23").Value = 2000000
$ws.Range("M23").Value = -22079.572
$ws.Range("N23").Value = -2000460

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 15000
$ws.Range("K37").Value = 15000
$ws.Range("M37").Value = -14797
$ws.Range("H46").Value = 93806.54
$ws.Range("J46").Value = 93806.54
$ws.Range("L46").Value = 93806.54
$ws.Range("N46").Value = -94268.54
$ws.Range("H62").Value = 5796.357
$ws.Range("I62").Value = 5727
$ws.Range("J62").Value = 5848.375
$ws.Range("K62").Value = 5727
$ws.Range("L62").Value = 5848.375
$ws.Range("M62").Value = -5103
$ws.Range("N62").Value = -7096.375
$ws.Range("H65").Value = 5796.357
$ws.Range("I65").Value = 5727
$ws.Range("J65").Value = 5848.375
$ws.Range("K65").Value = 28635
$ws.Range("L65").Value = 29241.875
$ws.Range("M65").Value = -25515
$ws.Range("N65").Value = -35481.875
$ws.Range("H92").Value = 84998.25
$ws.Range("J92").Value = 84998.25
$ws.Range("L92").Value = 84998.25
$ws.Range("N92").Value = -89990.25
$ws.Range("H107").Value = 71428980
$ws.Range("I107").Value = 458.75
$ws.Range("K107").Value = 1376.25
$ws.Range("M107").Value = 543.75
$ws.Range("H123").Value = 87509.71
$ws.Range("J123").Value = 87509.71
$ws.Range("L123").Value = 87509.71
$ws.Range("N123").Value = -97309.71
$ws.Range("H125").Value = 85244.625
$ws.Range("J125").Value = 85244.625
$ws.Range("L125").Value = 85244.625
$ws.Range("N125").Value = -95084.625
$ws.Range("H126").Value = 2421.05
$ws.Range("I126").Value = 2215.9412
$ws.Range("K126").Value = 6647.823600000001
$ws.Range("M126").Value = -4177.823600000001
$ws.Range("H132").Value = 32968.312
$ws.Range("I132").Value = 43200.082
$ws.Range("K132").Value = 129600.246
$ws.Range("M132").Value = -127070.246
$ws.Range("H134").Value = 93806.54
$ws.Range("J134").Value = 93806.54
$ws.Range("L134").Value = 281419.62
$ws.Range("N134").Value = -286489.62
$ws.Range("H140").Value = 73988.4
$ws.Range("J140").Value = 73988.4
$ws.Range("L140").Value = 73988.4
$ws.Range("N140").Value = -84348.4
$ws.Range("H141").Value = 109222.75
$ws.Range("J141").Value = 109222.75
$ws.Range("L141").Value = 109222.75
$ws.Range("N141").Value = -119582.75
